$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows starting at row 24 (pushes "Sector Distribution Details"
# and everything below it down to make room for the new enterprise-size
# breakdown table).
$ws.Rows("24:29").Insert()

# New mini-table: number of employees / assets / turnover by enterprise size
# class, mirroring the "title" (bold) style used by the other section header
# rows (e.g. B11:D11) and the plain "Normal" style used by their data rows.
$ws.Range("B23:D23").Style = "title"
$ws.Range("B23").Value = "Number of employees"
$ws.Range("C23").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D23").Value = "Turnover (local currency, unless noted otherwise)"

$ws.Range("A24:D27").Style = "Normal"

$ws.Range("A24").Value = "Micro"

$ws.Range("A25").Value = "Small"
$ws.Range("B25").Value = "1-99"

$ws.Range("A26").Value = "Medium"
$ws.Range("B26").Value = "100-499"

$ws.Range("A27").Value = "Large"
$ws.Range("B27").Value = ">=500"

# The hyperlinked source URL moved from A48 to A54 along with the row
# insert; re-anchor the hyperlink so it still points at the right cell.
$ws.Range("A48").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A54"), "http://www.ic.gc.ca/eic/site/061.nsf/eng/h_02800.html")
